# Edit: 2020-06-18 commit
#
# 1) Slide 5's table (graphicFrame, Shapes.Item(2)) gets a new built-in
#    table style GUID.
# 2) The deck's theme color scheme (the "Red Violet"/"Integral" palette
#    used by the slide master / all slides) is swapped back to the
#    stock Office blue palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{76452B4E-F0D6-431D-AEE8-932AE73CFC65}")

# --- 2. Theme colour scheme ------------------------------------------------
# Colour order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. RGB() isn't available in this host, so the
# packed 0xBBGGRR integers are spelled out explicitly (Office theme colours).
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

$themeColors.Item(1).RGB  = 0        # dk1      000000
$themeColors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388  # dk2      44546A
$themeColors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501  # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407    # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308 # accent5  4472C4
$themeColors.Item(10).RGB = 4697456  # accent6  70AD47
$themeColors.Item(11).RGB = 12673797 # hlink    0563C1
$themeColors.Item(12).RGB = 7491477  # folHlink 954F72
